# Evaluations: #3, finished analyzing video 4
#
# Fills in the VidID (A) / Class (B) columns for the "o1vRJ07KDZ0" / "Win"
# event block (rows 107-118, which already had these values on row 106 but
# were missing them on the following rows) and appends two new events
# (rows 119-120) for that same video ("Android Device Chooser" window
# config events).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$vidId = "o1vRJ07KDZ0"
$class = "Win"

# Rows 107-118 already have Focus/Action/min/sec/Notes - just backfill the
# VidID / Class columns that carry over from row 106.
for ($r = 107; $r -le 118; $r++) {
    $ws.Cells.Item($r, 1).Value = $vidId
    $ws.Cells.Item($r, 2).Value = $class
}

# New rows: two "Android Device Chooser" window-config/open events.
$ws.Range("A119").Value = $vidId
$ws.Range("B119").Value = $class
$ws.Range("C119").Value = "windowconfig"
$ws.Range("D119").Value = "open"
$ws.Range("E119").Value = 21
$ws.Range("F119").Value = 15
$ws.Range("G119").Value = "Android Device Chooser"

$ws.Range("A120").Value = $vidId
$ws.Range("B120").Value = $class
$ws.Range("C120").Value = "windowconfig"
$ws.Range("D120").Value = "open"
$ws.Range("E120").Value = 22
$ws.Range("F120").Value = 14
$ws.Range("G120").Value = "Android Device Chooser"

# Match the author's final selection/scroll state recorded in the diff.
$ws.Range("A106:B120").Select() | Out-Null
